# Apply "1st changes of mifos to finflux":
# Insert a new blank column before column N (14) on the "Repayment schedule" sheet,
# shifting the old N/O/P columns (Outstanding / Date / Disbursement) one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new column before column N - this shifts N,O,P -> O,P,Q
$ws.Columns.Item(14).Insert()

# New column N should look like column M (same width, no bestFit flag)
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Move the selection the way it ended up after the edit
$ws.Range("R6").Select() | Out-Null
